$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price reading was inserted at the top of the data (row 8),
# pushing the existing rows 8-14 down to rows 9-15.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with the latest reading.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C8").Value = 'Ñuble'
$ws.Range("D8").Value = 45167
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100102
$ws.Range("H8").Value = 'Cítricos'
$ws.Range("I8").Value = 100102006
$ws.Range("J8").Value = 'Pomelo'
$ws.Range("K8").Value = 'Start Ruby'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range("Q8").Value = '$/caja 14 kilos empedrada'
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1143
$ws.Range("T8").Value = 14
